# Apply "Deleted some columns" edit to DB_Table_Structure.xlsx
# Two stacked mini field-lists live in the sheet: one in column A (rows 1-48)
# and one in column D (rows 1-48). Several rows were removed from each list
# (YEAR, QUARTER, TAIL_NUM, WHEELS_OFF, WHEELS_ON from the column-A list;
# CITY_NAME, STATE_ABR from the column-D list), so every remaining entry
# shifts up within its own list and the now-unused rows at the bottom
# become empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row codes:
#  H = section header (bold font, boxed border)      -> style like s="2"
#  D = normal data row (regular font, boxed border)  -> style like s="3"
#  B = blank row that stays present (no text/border) -> style like s="1"
#  X = row no longer used at all -> cell removed entirely

$A_DATA = @(
    @("H", "FLIGHT_SCHED"),
    @("D", "FLIGHT_ID"),
    @("D", "MKT_UNIQUE_CARRIER"),
    @("D", "MKT_CARRIER_FL_NUM"),
    @("D", "FL_DATE"),
    @("D", "MONTH"),
    @("D", "DAY_OF_MONTH"),
    @("D", "DAY_OF_WEEK"),
    @("D", "ORIGIN"),
    @("D", "DEST"),
    @("D", "CRS_DEP_TIME"),
    @("D", "DEP_TIME_BLK"),
    @("D", "CRS_ARR_TIME"),
    @("D", "ARR_TIME_BLK"),
    @("D", "CRS_ELAPSED_TIME"),
    @("D", "DISTANCE"),
    @("D", "DISTANCE_GROUP"),
    @("B", $null),
    @("H", "FLIGHT_ACTUALS"),
    @("D", "FLIGHT_ID"),
    @("D", "DEP_TIME"),
    @("D", "TAXI_OUT"),
    @("D", "AIR_TIME"),
    @("D", "TAXI_IN"),
    @("D", "ARR_TIME"),
    @("D", "ACTUAL_ELAPSED_TIME"),
    @("B", $null),
    @("H", "DELAYS_CANCELS"),
    @("D", "FLIGHT_ID"),
    @("D", "DEP_DEL15"),
    @("D", "DEP_DELAY"),
    @("D", "DEP_DELAY_GROUP"),
    @("D", "ARR_DEL15"),
    @("D", "ARR_DELAY"),
    @("D", "ARR_DELAY_GROUP"),
    @("D", "TOTAL_DELAY"),
    @("D", "CARRIER_DELAY"),
    @("D", "WEATHER_DELAY"),
    @("D", "NAS_DELAY"),
    @("D", "SECURITY_DELAY"),
    @("D", "LATE_AIRCRAFT"),
    @("D", "CANCELLED"),
    @("D", "CANCELLATION_CODE"),
    @("X", $null),
    @("X", $null),
    @("X", $null),
    @("X", $null),
    @("X", $null)
)

$D_DATA = @(
    @("H", "AIRLINE"),
    @("D", "MKT_UNIQUE_CARRIER"),
    @("D", "CARRIER_NAME"),
    @("B", $null),
    @("H", "LOCATIONS"),
    @("D", "AIRPORT_CODE"),
    @("D", "STATE_NM"),
    @("B", $null),
    @("H", "CANCELLATION"),
    @("D", "CANCELLATION_CODE"),
    @("D", "CANCELATION_REASON"),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("B", $null),
    @("X", $null),
    @("X", $null)
)

function Set-FieldCell {
    param($ColumnLetter, $RowNumber, $Code, $Text)

    $cell = $ws.Range("$ColumnLetter$RowNumber")

    if ($Code -eq "X") {
        # Entirely empty the cell - no value, no formatting left behind.
        $cell.Clear()
        return
    }

    if ($Code -eq "B") {
        $cell.ClearContents()
        $cell.Font.Bold = $false
        $cell.Borders.LineStyle = 0
        return
    }

    $cell.Value = $Text
    if ($Code -eq "H") {
        $cell.Font.Bold = $true
    } else {
        $cell.Font.Bold = $false
    }
    $cell.Borders.LineStyle = 1
}

for ($i = 0; $i -lt $A_DATA.Count; $i++) {
    $row = $i + 1
    $entry = $A_DATA[$i]
    Set-FieldCell "A" $row $entry[0] $entry[1]
}

for ($i = 0; $i -lt $D_DATA.Count; $i++) {
    $row = $i + 1
    $entry = $D_DATA[$i]
    Set-FieldCell "D" $row $entry[0] $entry[1]
}

# Selection moved to D26 in the saved file.
$ws.Range("D26").Select()
